$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-06-05 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-06-06 Friday", 2) | Out-Null
$d.Content.Find.Execute("54+26=80", $true, $false, $false, $false, $false, $true, 1, $false, "14-6=8", 2) | Out-Null
$d.Content.Find.Execute("14+61=75", $true, $false, $false, $false, $false, $true, 1, $false, "71+24=95", 2) | Out-Null
$d.Content.Find.Execute("21+52=73", $true, $false, $false, $false, $false, $true, 1, $false, "7+53=60", 2) | Out-Null
$d.Content.Find.Execute("58-47=11", $true, $false, $false, $false, $false, $true, 1, $false, "80-3=77", 2) | Out-Null
$d.Content.Find.Execute("80-66=14", $true, $false, $false, $false, $false, $true, 1, $false, "92-34=58", 2) | Out-Null
$d.Content.Find.Execute("6+69=75", $true, $false, $false, $false, $false, $true, 1, $false, "80-56=24", 2) | Out-Null
$d.Content.Find.Execute("43+39=82", $true, $false, $false, $false, $false, $true, 1, $false, "44-15=29", 2) | Out-Null
$d.Content.Find.Execute("68+7=75", $true, $false, $false, $false, $false, $true, 1, $false, "28+20=48", 2) | Out-Null
$d.Content.Find.Execute("95-78=17", $true, $false, $false, $false, $false, $true, 1, $false, "81-70=11", 2) | Out-Null
$d.Content.Find.Execute("1+16=17", $true, $false, $false, $false, $false, $true, 1, $false, "21+18=39", 2) | Out-Null
$d.Content.Find.Execute("88-79=9", $true, $false, $false, $false, $false, $true, 1, $false, "11+86=97", 2) | Out-Null
$d.Content.Find.Execute("46+40=86", $true, $false, $false, $false, $false, $true, 1, $false, "67-38=29", 2) | Out-Null
$d.Content.Find.Execute("80+16=96", $true, $false, $false, $false, $false, $true, 1, $false, "49+43=92", 2) | Out-Null
$d.Content.Find.Execute("53+46=99", $true, $false, $false, $false, $false, $true, 1, $false, "69-28=41", 2) | Out-Null
$d.Content.Find.Execute("38+23=61", $true, $false, $false, $false, $false, $true, 1, $false, "89-9=80", 2) | Out-Null
$d.Content.Find.Execute("42-17=25", $true, $false, $false, $false, $false, $true, 1, $false, "71-31=40", 2) | Out-Null
$d.Content.Find.Execute("37+21=58", $true, $false, $false, $false, $false, $true, 1, $false, "52-23=29", 2) | Out-Null
$d.Content.Find.Execute("78-57=21", $true, $false, $false, $false, $false, $true, 1, $false, "12-1=11", 2) | Out-Null
$d.Content.Find.Execute("89-0=89", $true, $false, $false, $false, $false, $true, 1, $false, "63-61=2", 2) | Out-Null
$d.Content.Find.Execute("98+1=99", $true, $false, $false, $false, $false, $true, 1, $false, "78-8=70", 2) | Out-Null
$d.Content.Find.Execute("9-1=8", $true, $false, $false, $false, $false, $true, 1, $false, "0+56=56", 2) | Out-Null
$d.Content.Find.Execute("74-62=12", $true, $false, $false, $false, $false, $true, 1, $false, "52-4=48", 2) | Out-Null
$d.Content.Find.Execute("8+53=61", $true, $false, $false, $false, $false, $true, 1, $false, "54+40=94", 2) | Out-Null
$d.Content.Find.Execute("87-34=53", $true, $false, $false, $false, $false, $true, 1, $false, "16-14=2", 2) | Out-Null
$d.Content.Find.Execute("18+32=50", $true, $false, $false, $false, $false, $true, 1, $false, "22-8=14", 2) | Out-Null
$d.Content.Find.Execute("58-37=21", $true, $false, $false, $false, $false, $true, 1, $false, "2+38=40", 2) | Out-Null
$d.Content.Find.Execute("83-16=67", $true, $false, $false, $false, $false, $true, 1, $false, "11+43=54", 2) | Out-Null
$d.Content.Find.Execute("58-53=5", $true, $false, $false, $false, $false, $true, 1, $false, "76-68=8", 2) | Out-Null
$d.Content.Find.Execute("99-48=51", $true, $false, $false, $false, $false, $true, 1, $false, "55+44=99", 2) | Out-Null
$d.Content.Find.Execute("49+34=83", $true, $false, $false, $false, $false, $true, 1, $false, "96-17=79", 2) | Out-Null
$d.Content.Find.Execute("98-14=84", $true, $false, $false, $false, $false, $true, 1, $false, "62+31=93", 2) | Out-Null
$d.Content.Find.Execute("25+24=49", $true, $false, $false, $false, $false, $true, 1, $false, "38-27=11", 2) | Out-Null
$d.Content.Find.Execute("33+13=46", $true, $false, $false, $false, $false, $true, 1, $false, "84-56=28", 2) | Out-Null
$d.Content.Find.Execute("57+17=74", $true, $false, $false, $false, $false, $true, 1, $false, "30-19=11", 2) | Out-Null
$d.Content.Find.Execute("66-26=40", $true, $false, $false, $false, $false, $true, 1, $false, "25+15=40", 2) | Out-Null
$d.Content.Find.Execute("22-16=6", $true, $false, $false, $false, $false, $true, 1, $false, "4+11=15", 2) | Out-Null
$d.Content.Find.Execute("66+24=90", $true, $false, $false, $false, $false, $true, 1, $false, "36+11=47", 2) | Out-Null
$d.Content.Find.Execute("7+66=73", $true, $false, $false, $false, $false, $true, 1, $false, "18+15=33", 2) | Out-Null
$d.Content.Find.Execute("28+17=45", $true, $false, $false, $false, $false, $true, 1, $false, "4+2=6", 2) | Out-Null
$d.Content.Find.Execute("83-70=13", $true, $false, $false, $false, $false, $true, 1, $false, "14+22=36", 2) | Out-Null
$d.Content.Find.Execute("18+65=83", $true, $false, $false, $false, $false, $true, 1, $false, "27+50=77", 2) | Out-Null
$d.Content.Find.Execute("89-54=35", $true, $false, $false, $false, $false, $true, 1, $false, "59+11=70", 2) | Out-Null
$d.Content.Find.Execute("3+88=91", $true, $false, $false, $false, $false, $true, 1, $false, "86+4=90", 2) | Out-Null
$d.Content.Find.Execute("11+72=83", $true, $false, $false, $false, $false, $true, 1, $false, "30+13=43", 2) | Out-Null
$d.Content.Find.Execute("45-25=20", $true, $false, $false, $false, $false, $true, 1, $false, "82-63=19", 2) | Out-Null
$d.Content.Find.Execute("41+12=53", $true, $false, $false, $false, $false, $true, 1, $false, "74+9=83", 2) | Out-Null
$d.Content.Find.Execute("24+27=51", $true, $false, $false, $false, $false, $true, 1, $false, "51+4=55", 2) | Out-Null
$d.Content.Find.Execute("95-74=21", $true, $false, $false, $false, $false, $true, 1, $false, "30+30=60", 2) | Out-Null
$d.Content.Find.Execute("6+47=53", $true, $false, $false, $false, $false, $true, 1, $false, "40+5=45", 2) | Out-Null
$d.Content.Find.Execute("91-62=29", $true, $false, $false, $false, $false, $true, 1, $false, "52-5=47", 2) | Out-Null
$d.Content.Find.Execute("39+43=82", $true, $false, $false, $false, $false, $true, 1, $false, "38-22=16", 2) | Out-Null
$d.Content.Find.Execute("60+33=93", $true, $false, $false, $false, $false, $true, 1, $false, "58-9=49", 2) | Out-Null
$d.Content.Find.Execute("0+83=83", $true, $false, $false, $false, $false, $true, 1, $false, "43+14=57", 2) | Out-Null
$d.Content.Find.Execute("71+13=84", $true, $false, $false, $false, $false, $true, 1, $false, "61-33=28", 2) | Out-Null
$d.Content.Find.Execute("50+15=65", $true, $false, $false, $false, $false, $true, 1, $false, "58+22=80", 2) | Out-Null
$d.Content.Find.Execute("71-26=45", $true, $false, $false, $false, $false, $true, 1, $false, "86-8=78", 2) | Out-Null
$d.Content.Find.Execute("37-25=12", $true, $false, $false, $false, $false, $true, 1, $false, "30+23=53", 2) | Out-Null
$d.Content.Find.Execute("96-8=88", $true, $false, $false, $false, $false, $true, 1, $false, "32+47=79", 2) | Out-Null
$d.Content.Find.Execute("27-14=13", $true, $false, $false, $false, $false, $true, 1, $false, "83-78=5", 2) | Out-Null
$d.Content.Find.Execute("40+23=63", $true, $false, $false, $false, $false, $true, 1, $false, "99-81=18", 2) | Out-Null
$d.Content.Find.Execute("64+1=65", $true, $false, $false, $false, $false, $true, 1, $false, "5+79=84", 2) | Out-Null
$d.Content.Find.Execute("8+68=76", $true, $false, $false, $false, $false, $true, 1, $false, "54-26=28", 2) | Out-Null
$d.Content.Find.Execute("29+59=88", $true, $false, $false, $false, $false, $true, 1, $false, "9+64=73", 2) | Out-Null
$d.Content.Find.Execute("53-42=11", $true, $false, $false, $false, $false, $true, 1, $false, "44+54=98", 2) | Out-Null
$d.Content.Find.Execute("59+3=62", $true, $false, $false, $false, $false, $true, 1, $false, "95-25=70", 2) | Out-Null
$d.Content.Find.Execute("32+37=69", $true, $false, $false, $false, $false, $true, 1, $false, "98-52=46", 2) | Out-Null
$d.Content.Find.Execute("66-21=45", $true, $false, $false, $false, $false, $true, 1, $false, "38+40=78", 2) | Out-Null
$d.Content.Find.Execute("71-46=25", $true, $false, $false, $false, $false, $true, 1, $false, "40+42=82", 2) | Out-Null
$d.Content.Find.Execute("13-9=4", $true, $false, $false, $false, $false, $true, 1, $false, "72-1=71", 2) | Out-Null
$d.Content.Find.Execute("17+72=89", $true, $false, $false, $false, $false, $true, 1, $false, "67+32=99", 2) | Out-Null
$d.Content.Find.Execute("43-27=16", $true, $false, $false, $false, $false, $true, 1, $false, "62-23=39", 2) | Out-Null
$d.Content.Find.Execute("57+26=83", $true, $false, $false, $false, $false, $true, 1, $false, "20+20=40", 2) | Out-Null
$d.Content.Find.Execute("17+33=50", $true, $false, $false, $false, $false, $true, 1, $false, "18+28=46", 2) | Out-Null
$d.Content.Find.Execute("9-7=2", $true, $false, $false, $false, $false, $true, 1, $false, "72-26=46", 2) | Out-Null
$d.Content.Find.Execute("79+9=88", $true, $false, $false, $false, $false, $true, 1, $false, "77-69=8", 2) | Out-Null
$d.Content.Find.Execute("23+21=44", $true, $false, $false, $false, $false, $true, 1, $false, "23+56=79", 2) | Out-Null
$d.Content.Find.Execute("66-64=2", $true, $false, $false, $false, $false, $true, 1, $false, "0+50=50", 2) | Out-Null
$d.Content.Find.Execute("88-86=2", $true, $false, $false, $false, $false, $true, 1, $false, "66-35=31", 2) | Out-Null
$d.Content.Find.Execute("93-10=83", $true, $false, $false, $false, $false, $true, 1, $false, "79-8=71", 2) | Out-Null
$d.Content.Find.Execute("58+6=64", $true, $false, $false, $false, $false, $true, 1, $false, "56+6=62", 2) | Out-Null
$d.Content.Find.Execute("25+63=88", $true, $false, $false, $false, $false, $true, 1, $false, "81-53=28", 2) | Out-Null
$d.Content.Find.Execute("36-10=26", $true, $false, $false, $false, $false, $true, 1, $false, "33-17=16", 2) | Out-Null
$d.Content.Find.Execute("5+26=31", $true, $false, $false, $false, $false, $true, 1, $false, "46+11=57", 2) | Out-Null
$d.Content.Find.Execute("62+24=86", $true, $false, $false, $false, $false, $true, 1, $false, "96-33=63", 2) | Out-Null
$d.Content.Find.Execute("44+11=55", $true, $false, $false, $false, $false, $true, 1, $false, "76-51=25", 2) | Out-Null
$d.Content.Find.Execute("45+3=48", $true, $false, $false, $false, $false, $true, 1, $false, "25+67=92", 2) | Out-Null
$d.Content.Find.Execute("26+66=92", $true, $false, $false, $false, $false, $true, 1, $false, "54+23=77", 2) | Out-Null
$d.Content.Find.Execute("29-26=3", $true, $false, $false, $false, $false, $true, 1, $false, "40-26=14", 2) | Out-Null
$d.Content.Find.Execute("15+68=83", $true, $false, $false, $false, $false, $true, 1, $false, "26-4=22", 2) | Out-Null
$d.Content.Find.Execute("78+20=98", $true, $false, $false, $false, $false, $true, 1, $false, "75-42=33", 2) | Out-Null
$d.Content.Find.Execute("5+61=66", $true, $false, $false, $false, $false, $true, 1, $false, "16+55=71", 2) | Out-Null
$d.Content.Find.Execute("54-5=49", $true, $false, $false, $false, $false, $true, 1, $false, "10+55=65", 2) | Out-Null
$d.Content.Find.Execute("64+11=75", $true, $false, $false, $false, $false, $true, 1, $false, "11+74=85", 2) | Out-Null
$d.Content.Find.Execute("83+0=83", $true, $false, $false, $false, $false, $true, 1, $false, "46-9=37", 2) | Out-Null
$d.Content.Find.Execute("87+6=93", $true, $false, $false, $false, $false, $true, 1, $false, "96-41=55", 2) | Out-Null
$d.Content.Find.Execute("99-28=71", $true, $false, $false, $false, $false, $true, 1, $false, "76-51=25", 2) | Out-Null
$d.Content.Find.Execute("34-0=34", $true, $false, $false, $false, $false, $true, 1, $false, "48+19=67", 2) | Out-Null
$d.Content.Find.Execute("70-22=48", $true, $false, $false, $false, $false, $true, 1, $false, "25+18=43", 2) | Out-Null
$d.Content.Find.Execute("47-4=43", $true, $false, $false, $false, $false, $true, 1, $false, "69-50=19", 2) | Out-Null
$d.Content.Find.Execute("19+62=81", $true, $false, $false, $false, $false, $true, 1, $false, "7-1=6", 2) | Out-Null
